$d = $word.ActiveDocument

# 1. Title: "Team Buttercup Minutes" -> "Team Buttercup Client Minutes 1"
$d.Content.Find.Execute("Team Buttercup Minutes", $true, $false, $false, $false, $false, $true, 1, $false, "Team Buttercup Client Minutes 1", 2)

# 2. Fix typo: "Minuets" -> "Minutes" (Meeting Client Minutes 1 | Date: ...)
$d.Content.Find.Execute("Minuets", $true, $false, $false, $false, $false, $true, 1, $false, "Minutes", 2)

# 3. Move the "_GoBack" bookmark from after "figure" (in the Next Meeting Agenda
#    Items bullet) to the empty paragraph right after the Date/Time/Location line.
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
}
$targetPara = $d.Paragraphs.Item(3)
$targetRange = $targetPara.Range
$targetRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $targetRange)
